$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" number-looking string need to be forced
# back to text (matching the source data, which stores these as plain strings)
# so Excel COM does not auto-convert them to numeric values / drop formatting
# such as trailing zeros. We briefly mark the cell as Text, set the value, then
# restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.609"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000291"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "708.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.935"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "105.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "596.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.352"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.132"
$ws.Range("D49").Style = "Normal"

$ws.Range("D2").Value = "72.373.42"
$ws.Range("E2").Value = "  +4.62%  "
$ws.Range("D3").Value = "3.624.03"
$ws.Range("E3").Value = "  +7.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("D7").Value = "3.615.45"
$ws.Range("E7").Value = "  +7.14%  "
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("E12").Value = "  +4.37%  "
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("D15").Value = "4.205.58"
$ws.Range("E15").Value = "  +7.12%  "
$ws.Range("E16").Value = "  +3.95%  "
$ws.Range("D17").Value = "72.411.40"
$ws.Range("E17").Value = "  +4.51%  "
$ws.Range("D18").Value = "3.596.26"
$ws.Range("E18").Value = "  +5.97%  "
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("E23").Value = "  +9.14%  "
$ws.Range("E24").Value = "  +4.37%  "
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("E27").Value = "  +5.13%  "
$ws.Range("E28").Value = "  +4.55%  "
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("E30").Value = "  +4.61%  "
$ws.Range("E31").Value = "  +6.91%  "
$ws.Range("E32").Value = "  +14.85%  "
$ws.Range("E33").Value = "  +7.49%  "
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "3.647.45"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").Value = "0.0₃0783"
$ws.Range("E40").Value = "  +9.08%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +7.12%  "
$ws.Range("E43").Value = "  +4.91%  "
$ws.Range("E44").Value = "  +6.66%  "
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("E47").Value = "  +5.38%  "
$ws.Range("E48").Value = "  +5.75%  "
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  +0.44%  "
